# Change "We aim to combine multiple types of substrate in part to address
# this question." to "...substrates..." -- splitting the sentence into three
# runs (identical formatting) around the inserted "s", matching how Word
# leaves a freshly-edited word as its own run rather than folding it back
# into one contiguous run.

$d = $word.ActiveDocument

# 1) Perform the actual text edit: substrate -> substrates
$text = $d.Content.Text
$idx = $text.IndexOf("substrate in part to address this question.")
$r = $d.Range($idx, $idx + 9)
$r.Text = "substrates"

# 2) Force "substrates" to live in its own run (Word does this naturally when
#    you type/correct a word -- nudge formatting off and back on to split it
#    from its neighbours without any visible change).
$rWord = $d.Range($idx, $idx + 10)
$rWord.Bold = 1
$rWord.Bold = 0

# 3) Force the rest of the sentence ("We aim to combine multiple types of " /
#    " in part to address this question. ") to split away from the
#    preceding, unrelated sentence ("...fuller nutrient profile. ") the same
#    way.
$text2 = $d.Content.Text
$idxSentence = $text2.IndexOf("We aim to combine multiple types of substrates")
$sentence = "We aim to combine multiple types of substrates in part to address this question. "
$rSentence = $d.Range($idxSentence, $idxSentence + $sentence.Length)
$rSentence.Bold = 1
$rSentence.Bold = 0
